# Insert a new variable row ("climate_change_factor_gnrl_hydropower_availability")
# above the existing "elasticity_gnrl_rate_occupancy_to_gdppc" row on the first
# sheet, pushing all subsequent "General" rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 4 (shifts old rows 4-11 down to 5-12).
$ws.Rows("4:4").Insert()

# Populate the newly inserted row with the new variable's metadata/values.
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5
$ws.Range("J4:AS4").Value = 1
